$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"

$ws.Range("B1").Select()
